$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows for "Бокситогорский МР" (this & next input data), years 2016 down to 2011
$name = "Бокситогорский МР"

$data = @(
    @(2016, 50.756, 9.8539999999999992, 497.40879999999999, 33918.400000000001, 30.7, 894.31299999999999, 22261.996999999999, 196.75200000000001, 5, 1640.018, 14.432, 59.8, 33.5, 169),
    @(2015, 51.298000000000002, 9.8369999999999997, 559.14819999999997, 31340.5, 36.43, 1381.8050000000001, 19849.96, 159.67099999999999, 8.1, 1348.8230000000001, 17.542000000000002, 59.8, 33.5, 23),
    @(2014, 51.941000000000003, 9.9879999999999995, 597.32150000000001, 28366.2, 36.43, 814.22699999999998, 17254.972000000002, 173.51900000000001, 29.9, 967.63780000000008, 16.149000000000001, 59.8, 33.5, -257),
    @(2013, 52.34, 10.254, 711.82400000000007, 25409.8, 36.43, 691.16099999999994, 17653.433000000001, 85.984999999999999, 9.9, 821.76700000000005, 24.484999999999999, 59.8, 33.5, 161),
    @(2012, 52.956000000000003, 10.59, 751.97519999999997, 23164, 36.43, 955.6, 16316.396500000001, 230.62899999999999, 5.4, 714.0385, 32.409999999999997, 59.8, 33.5, -76),
    @(2011, 53.338000000000001, 11.303000000000001, 1002.7544, 20279, 36.43, 577.28300000000002, 15980.2601, 309.06599999999997, 5.6, 633.4, 29.247, 59.8, 33.5, -150)
)

$row = 15
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $r[0]
    $ws.Cells.Item($row, 3).Value = $r[1]
    $ws.Cells.Item($row, 4).Value = $r[2]
    $ws.Cells.Item($row, 5).Value = $r[3]
    $ws.Cells.Item($row, 6).Value = $r[4]
    $ws.Cells.Item($row, 7).Value = $r[5]
    $ws.Cells.Item($row, 8).Value = $r[6]
    $ws.Cells.Item($row, 9).Value = $r[7]
    $ws.Cells.Item($row, 10).Value = $r[8]
    $ws.Cells.Item($row, 11).Value = $r[9]
    $ws.Cells.Item($row, 12).Value = $r[10]
    $ws.Cells.Item($row, 13).Value = $r[11]
    $ws.Cells.Item($row, 14).Value = $r[12]
    $ws.Cells.Item($row, 15).Value = $r[13]
    $ws.Cells.Item($row, 16).Value = $r[14]

    $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 2)).HorizontalAlignment = -4108
    $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 2)).Font.Bold = $true
    $ws.Range($ws.Cells.Item($row, 3), $ws.Cells.Item($row, 16)).HorizontalAlignment = -4108

    $row++
}

$ws.Range("M30").Select()
